$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I1 ("I0") and J1 ("IF"), styled like the
# existing header cells (e.g. H1) by copying H1's formatting over.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data row: new values for columns I and J
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
